$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update harvester column (B) and add experimentDesign column (D) for data rows 2-19
for ($r = 2; $r -le 19; $r++) {
    $ws.Cells.Item($r, 2).Value = "S.GISH"
    $ws.Cells.Item($r, 4).Value = "90minuteInduction"
}

# Update selection to match the diff (active cell D3, selection D3:D19)
$ws.Range("D3:D19").Select()
